$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1: the scraper re-ran and a batch of matches came back with their
# "home/away" pairing flipped relative to the previous sheet order. For each
# of the following adjacent row pairs, the match-detail columns F:V (teams,
# scores, odds, timestamps, url) swap places while the leading index/meta
# columns A:E (row index, pais/torneio/temporada, data_partida) stay put.
# ---------------------------------------------------------------------------
$swapPairs = @(
    @(19, 20),
    @(22, 23),
    @(37, 38),
    @(62, 63),
    @(74, 75),
    @(84, 85),
    @(96, 97),
    @(98, 99),
    @(106, 107),
    @(110, 111),
    @(133, 134),
    @(167, 168),
    @(176, 177)
)

foreach ($pair in $swapPairs) {
    $a = $pair[0]
    $b = $pair[1]
    $rangeA = $ws.Range("F$($a):V$($a)")
    $rangeB = $ws.Range("F$($b):V$($b)")
    $valA = $rangeA.Value()
    $valB = $rangeB.Value()
    $rangeA.Value = $valB
    $rangeB.Value = $valA
}

# ---------------------------------------------------------------------------
# Part 2: four new matches appended at the bottom of the sheet (rows 193-196),
# extending the used range from A1:V192 to A1:V196.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=193; A=192; E=45257.95833333334; F="Estudiantes L.P."; G=1; H="Lanus"; I=1;
       J=1.83; K="12/11/2023 21:12"; L=1.72; M="27/11/2023 22:58";
       N=3.33; O="12/11/2023 21:12"; P=3.69; Q="27/11/2023 22:58";
       R=5.02; S="12/11/2023 21:12"; T=5.37; U="27/11/2023 22:59";
       V="https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/estudiantes-l-p-lanus/hEVbm0yn/" },
    @{ Row=194; A=193; E=45257.95833333334; F="San Lorenzo"; G=2; H="Central Cordoba"; I=0;
       J=1.92; K="12/11/2023 21:12"; L=2.08; M="27/11/2023 22:58";
       N=3.08; O="12/11/2023 21:12"; P=2.96; Q="27/11/2023 22:58";
       R=5.02; S="12/11/2023 21:12"; T=4.6; U="27/11/2023 22:58";
       V="https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/san-lorenzo-central-cordoba-santiago-del-estero/zwL3nKjh/" },
    @{ Row=195; A=194; E=45258.0625; F="Racing Club"; G=4; H="Belgrano"; I=1;
       J=1.92; K="21/11/2023 01:42"; L=2.06; M="28/11/2023 01:22";
       N=3.34; O="21/11/2023 01:42"; P=3.21; Q="28/11/2023 01:21";
       R=4.13; S="21/11/2023 01:42"; T=4.14; U="28/11/2023 00:58";
       V="https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/racing-club-ca-belgrano-de-cordoba/vBRfltMu/" },
    @{ Row=196; A=195; E=45258.0625; F="Newells Old Boys"; G=3; H="Defensa y Justicia"; I=0;
       J=1.78; K="12/11/2023 21:12"; L=1.69; M="28/11/2023 01:28";
       N=3.5; O="12/11/2023 21:12"; P=3.64; Q="28/11/2023 01:28";
       R=5; S="12/11/2023 21:12"; T=5.85; U="28/11/2023 01:28";
       V="https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/newells-old-boys-defensa-y-justicia/ETJ7ov6b/" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Bring over the formatting (bold/bordered index style on col A, the
    # datetime number format on col E) from the last existing data row
    # before filling in the values, so the new rows carry the same per-
    # column styles as every other row in the sheet.
    $ws.Range("A192:V192").Copy()
    $ws.Range("A$($rowNum):V$($rowNum)").PasteSpecial(-4122)

    $ws.Range("A$rowNum").Value = $r.A

    # B/C/D are identical on every row of this sheet; copy them verbatim
    # (as values) from row 192 instead of re-typing literals so the
    # "2023" season text in column D is not auto-coerced into a number.
    $ws.Range("B192:D192").Copy()
    $ws.Range("B$($rowNum):D$($rowNum)").PasteSpecial(-4163)

    $ws.Range("E$rowNum").Value = $r.E
    $ws.Range("F$rowNum").Value = $r.F
    $ws.Range("G$rowNum").Value = $r.G
    $ws.Range("H$rowNum").Value = $r.H
    $ws.Range("I$rowNum").Value = $r.I
    $ws.Range("J$rowNum").Value = $r.J
    $ws.Range("K$rowNum").Value = $r.K
    $ws.Range("L$rowNum").Value = $r.L
    $ws.Range("M$rowNum").Value = $r.M
    $ws.Range("N$rowNum").Value = $r.N
    $ws.Range("O$rowNum").Value = $r.O
    $ws.Range("P$rowNum").Value = $r.P
    $ws.Range("Q$rowNum").Value = $r.Q
    $ws.Range("R$rowNum").Value = $r.R
    $ws.Range("S$rowNum").Value = $r.S
    $ws.Range("T$rowNum").Value = $r.T
    $ws.Range("U$rowNum").Value = $r.U
    $ws.Range("V$rowNum").Value = $r.V
}

Write-Output "done"
